$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.141.50"
$ws.Range('E2').Value = "'  +0.09%  "
$ws.Range('D3').Value = "'2.570.51"
$ws.Range('E3').Value = "'  -1.05%  "
$ws.Range('E4').Value = "'  +0.38%  "
$ws.Range('D5').Value = "'505.41"
$ws.Range('E5').Value = "'  -0.73%  "
$ws.Range('D6').Value = "'152.06"
$ws.Range('E6').Value = "'  -3.42%  "
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('D8').Value = "'0.578"
$ws.Range('E8').Value = "'  -5.43%  "
$ws.Range('D9').Value = "'2.574.89"
$ws.Range('E9').Value = "'  -0.17%  "
$ws.Range('D10').Value = "'6.57"
$ws.Range('E10').Value = "'  +7.32%  "
$ws.Range('E11').Value = "'  +0.20%  "
$ws.Range('E12').Value = "'  +1.19%  "
$ws.Range('E13').Value = "'  +1.18%  "
$ws.Range('D14').Value = "'3.022.07"
$ws.Range('E14').Value = "'  +0.12%  "
$ws.Range('D15').Value = "'60.200.68"
$ws.Range('E15').Value = "'  +1.08%  "
$ws.Range('E16').Value = "'  -2.02%  "
$ws.Range('E17').Value = "'  +1.20%  "
$ws.Range('D18').Value = "'2.571.95"
$ws.Range('E18').Value = "'  -0.21%  "
$ws.Range('E19').Value = "'  -0.22%  "
$ws.Range('D20').Value = "'344.45"
$ws.Range('E20').Value = "'  +1.28%  "
$ws.Range('D21').Value = "'10.38"
$ws.Range('E21').Value = "'  -0.33%  "
$ws.Range('D22').Value = "'6.08"
$ws.Range('E22').Value = "'  +0.51%  "
$ws.Range('E23').Value = "'  -0.38%  "
$ws.Range('D24').Value = "'59.70"
$ws.Range('E24').Value = "'  -0.74%  "
$ws.Range('E25').Value = "'  -0.20%  "
$ws.Range('E26').Value = "'  +0.42%  "
$ws.Range('D27').Value = "'0.999"
$ws.Range('E27').Value = "'  -0.16%  "
$ws.Range('D28').Value = "'0.0₃0838"
$ws.Range('E28').Value = "'  -0.03%  "
$ws.Range('E29').Value = "'  +0.46%  "
$ws.Range('D31').Value = "'19.27"
$ws.Range('E31').Value = "'  -0.84%  "
$ws.Range('D32').Value = "'153.18"
$ws.Range('E32').Value = "'  -2.26%  "
$ws.Range('E33').Value = "'  -1.14%  "
$ws.Range('E34').Value = "'  +2.66%  "
$ws.Range('E35').Value = "'  +1.64%  "
$ws.Range('E36').Value = "'  -1.09%  "
$ws.Range('D37').Value = "'0.845"
$ws.Range('E37').Value = "'  +7.87%  "
$ws.Range('D38').Value = "'0.847"
$ws.Range('E38').Value = "'  -2.16%  "
$ws.Range('E39').Value = "'  +1.30%  "
$ws.Range('D40').Value = "'36.08"
$ws.Range('E40').Value = "'  +2.41%  "
$ws.Range('E41').Value = "'  -0.35%  "
$ws.Range('D42').Value = "'293.64"
$ws.Range('E42').Value = "'  -4.13%  "
$ws.Range('D43').Value = "'0.616"
$ws.Range('E43').Value = "'  -2.24%  "
$ws.Range('D44').Value = "'0.0991"
$ws.Range('E44').Value = "'  -2.71%  "
$ws.Range('E45').Value = "'  +0.70%  "
$ws.Range('D46').Value = "'0.0555"
$ws.Range('E46').Value = "'  -2.94%  "
$ws.Range('D47').Value = "'19.65"
$ws.Range('E47').Value = "'  +1.55%  "
$ws.Range('E48').Value = "'  -2.22%  "
$ws.Range('E49').Value = "'  -2.19%  "
$ws.Range('E50').Value = "'  +0.40%  "
$ws.Range('D51').Value = "'1.990.59"
$ws.Range('E51').Value = "'  +0.16%  "
